$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new column before column G (shifts G->H, H->I, I->J, J->K)
$ws.Columns("G").Insert()

# New header
$ws.Range("G1").Value = "py_plot_ready"

# Fill new column for rows 2-9 with "yes"
$ws.Range("G2:G9").Value = "yes"

# Set the new column width to fit header text
$ws.Columns("G").ColumnWidth = 14.109375

$ws.Range("G2").Select()
